$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - swap/reassign labels
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"

# Update row 4 values
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

# Update row 5 values
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1

# Update row 6 values
$ws.Range("B6").Value = 1
$ws.Range("E6").Value = 0

# Update row 7 values
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = 0

$wb.Save()
